$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 7).Value = 1.823996
$ws.Cells.Item(2, 8).Value = 3.647992
$ws.Cells.Item(2, 9).Value = 0.04519532258275597
$ws.Cells.Item(2, 10).Value = 0.03907135622362176
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 13).Value = 2.552614
$ws.Cells.Item(2, 14).Value = 5.105228
$ws.Cells.Item(2, 15).Value = 0.05395552785558979
$ws.Cells.Item(2, 16).Value = 0.04209357928847631
$ws.Cells.Item(2, 17).Value = 4.655957725544
$ws.Cells.Item(2, 18).Value = 18.623830902176
$ws.Cells.Item(2, 19).Value = 0.002438537486556256
$ws.Cells.Item(2, 20).Value = 0.001644653231107325

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 7).Value = 1.823996
$ws.Cells.Item(3, 8).Value = 3.647992
$ws.Cells.Item(3, 9).Value = 0.04519532258275597
$ws.Cells.Item(3, 10).Value = 0.03907135622362176
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 22.74715733333333
$ws.Cells.Item(3, 14).Value = 68.241472
$ws.Cells.Item(3, 15).Value = 0.4808149140975291
$ws.Cells.Item(3, 16).Value = 0.5626639618043182
$ws.Cells.Item(3, 17).Value = 41.49072398737066
$ws.Cells.Item(3, 18).Value = 248.944343924224
$ws.Cells.Item(3, 19).Value = 0.02173058514523793
$ws.Cells.Item(3, 20).Value = 0.02198404408585082

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 7).Value = 1.823996
$ws.Cells.Item(4, 8).Value = 3.647992
$ws.Cells.Item(4, 9).Value = 0.04519532258275597
$ws.Cells.Item(4, 10).Value = 0.03907135622362176
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.281676
$ws.Cells.Item(4, 14).Value = 0.845028
$ws.Cells.Item(4, 15).Value = 0.005953887765346076
$ws.Cells.Item(4, 16).Value = 0.006967417149436334
$ws.Cells.Item(4, 17).Value = 0.5137758972959999
$ws.Cells.Item(4, 18).Value = 3.082655383776
$ws.Cells.Item(4, 19).Value = 0.00026908787817634
$ws.Cells.Item(4, 20).Value = 0.0002722264374041983

$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 7).Value = 1.823996
$ws.Cells.Item(5, 8).Value = 3.647992
$ws.Cells.Item(5, 9).Value = 0.04519532258275597
$ws.Cells.Item(5, 10).Value = 0.03907135622362176
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.454891
$ws.Cells.Item(5, 14).Value = 1.364673
$ws.Cells.Item(5, 15).Value = 0.009615196038945605
$ws.Cells.Item(5, 16).Value = 0.01125198935842686
$ws.Cells.Item(5, 17).Value = 0.8297193644359999
$ws.Cells.Item(5, 18).Value = 4.978316186616
$ws.Cells.Item(5, 19).Value = 0.0004345618866765841
$ws.Cells.Item(5, 20).Value = 0.000439630484447497

$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 7).Value = 1.823996
$ws.Cells.Item(6, 8).Value = 3.647992
$ws.Cells.Item(6, 9).Value = 0.04519532258275597
$ws.Cells.Item(6, 10).Value = 0.03907135622362176
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 3.179912
$ws.Cells.Item(6, 14).Value = 9.539736
$ws.Cells.Item(6, 15).Value = 0.06721495317910355
$ws.Cells.Item(6, 16).Value = 0.07865694415746599
$ws.Cells.Item(6, 17).Value = 5.800146768352
$ws.Cells.Item(6, 18).Value = 34.800880610112
$ws.Cells.Item(6, 19).Value = 0.003037801491314424
$ws.Cells.Item(6, 20).Value = 0.003073233484637878

$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 7).Value = 1.823996
$ws.Cells.Item(7, 8).Value = 3.647992
$ws.Cells.Item(7, 9).Value = 0.04519532258275597
$ws.Cells.Item(7, 10).Value = 0.03907135622362176
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 13).Value = 18.0933415
$ws.Cells.Item(7, 14).Value = 36.186683
$ws.Cells.Item(7, 15).Value = 0.3824455210634858
$ws.Cells.Item(7, 16).Value = 0.2983661082418763
$ws.Cells.Item(7, 17).Value = 33.002182522634
$ws.Cells.Item(7, 18).Value = 132.008730090536
$ws.Cells.Item(7, 19).Value = 0.01728474869479444
$ws.Cells.Item(7, 20).Value = 0.01165756850017404

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 10.26769266666667
$ws.Cells.Item(8, 8).Value = 30.803078
$ws.Cells.Item(8, 9).Value = 0.2544148574068134
$ws.Cells.Item(8, 10).Value = 0.3299124650827103
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 13).Value = 2.552614
$ws.Cells.Item(8, 14).Value = 5.105228
$ws.Cells.Item(8, 15).Value = 0.05395552785558979
$ws.Cells.Item(8, 16).Value = 0.04209357928847631
$ws.Cells.Item(8, 17).Value = 26.20945604863067
$ws.Cells.Item(8, 18).Value = 157.256736291784
$ws.Cells.Item(8, 19).Value = 0.01372708792568922
$ws.Cells.Item(8, 20).Value = 0.01388719650721574

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 10.26769266666667
$ws.Cells.Item(9, 8).Value = 30.803078
$ws.Cells.Item(9, 9).Value = 0.2544148574068134
$ws.Cells.Item(9, 10).Value = 0.3299124650827103
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 22.74715733333333
$ws.Cells.Item(9, 14).Value = 68.241472
$ws.Cells.Item(9, 15).Value = 0.4808149140975291
$ws.Cells.Item(9, 16).Value = 0.5626639618043182
$ws.Cells.Item(9, 17).Value = 233.5608205389796
$ws.Cells.Item(9, 18).Value = 2102.047384850816
$ws.Cells.Item(9, 19).Value = 0.1223264578091921
$ws.Cells.Item(9, 20).Value = 0.1856298546520666

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 10.26769266666667
$ws.Cells.Item(10, 8).Value = 30.803078
$ws.Cells.Item(10, 9).Value = 0.2544148574068134
$ws.Cells.Item(10, 10).Value = 0.3299124650827103
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.281676
$ws.Cells.Item(10, 14).Value = 0.845028
$ws.Cells.Item(10, 15).Value = 0.005953887765346076
$ws.Cells.Item(10, 16).Value = 0.006967417149436334
$ws.Cells.Item(10, 17).Value = 2.892162599576
$ws.Cells.Item(10, 18).Value = 26.029463396184
$ws.Cells.Item(10, 19).Value = 0.001514757506836693
$ws.Cells.Item(10, 20).Value = 0.002298637767030091

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 10.26769266666667
$ws.Cells.Item(11, 8).Value = 30.803078
$ws.Cells.Item(11, 9).Value = 0.2544148574068134
$ws.Cells.Item(11, 10).Value = 0.3299124650827103
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.454891
$ws.Cells.Item(11, 14).Value = 1.364673
$ws.Cells.Item(11, 15).Value = 0.009615196038945605
$ws.Cells.Item(11, 16).Value = 0.01125198935842686
$ws.Cells.Item(11, 17).Value = 4.670680984832667
$ws.Cells.Item(11, 18).Value = 42.036128863494
$ws.Cells.Item(11, 19).Value = 0.002446248729186903
$ws.Cells.Item(11, 20).Value = 0.003712171546323028

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 10.26769266666667
$ws.Cells.Item(12, 8).Value = 30.803078
$ws.Cells.Item(12, 9).Value = 0.2544148574068134
$ws.Cells.Item(12, 10).Value = 0.3299124650827103
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 3.179912
$ws.Cells.Item(12, 14).Value = 9.539736
$ws.Cells.Item(12, 15).Value = 0.06721495317910355
$ws.Cells.Item(12, 16).Value = 0.07865694415746599
$ws.Cells.Item(12, 17).Value = 32.65035912304533
$ws.Cells.Item(12, 18).Value = 293.853232107408
$ws.Cells.Item(12, 19).Value = 0.01710048272866727
$ws.Cells.Item(12, 20).Value = 0.02594990634286269

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 10.26769266666667
$ws.Cells.Item(13, 8).Value = 30.803078
$ws.Cells.Item(13, 9).Value = 0.2544148574068134
$ws.Cells.Item(13, 10).Value = 0.3299124650827103
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 13).Value = 18.0933415
$ws.Cells.Item(13, 14).Value = 36.186683
$ws.Cells.Item(13, 15).Value = 0.3824455210634858
$ws.Cells.Item(13, 16).Value = 0.2983661082418763
$ws.Cells.Item(13, 17).Value = 185.7768698350457
$ws.Cells.Item(13, 18).Value = 1114.661219010274
$ws.Cells.Item(13, 19).Value = 0.09729982270724118
$ws.Cells.Item(13, 20).Value = 0.09843469826721217

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 0.279608
$ws.Cells.Item(14, 8).Value = 0.838824
$ws.Cells.Item(14, 9).Value = 0.0069281806301764
$ws.Cells.Item(14, 10).Value = 0.008984118197880725
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 13).Value = 2.552614
$ws.Cells.Item(14, 14).Value = 5.105228
$ws.Cells.Item(14, 15).Value = 0.05395552785558979
$ws.Cells.Item(14, 16).Value = 0.04209357928847631
$ws.Cells.Item(14, 17).Value = 0.7137312953120001
$ws.Cells.Item(14, 18).Value = 4.282387771872
$ws.Cells.Item(14, 19).Value = 0.0003738136429800404
$ws.Cells.Item(14, 20).Value = 0.0003781736916995352

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 0.279608
$ws.Cells.Item(15, 8).Value = 0.838824
$ws.Cells.Item(15, 9).Value = 0.0069281806301764
$ws.Cells.Item(15, 10).Value = 0.008984118197880725
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 22.74715733333333
$ws.Cells.Item(15, 14).Value = 68.241472
$ws.Cells.Item(15, 15).Value = 0.4808149140975291
$ws.Cells.Item(15, 16).Value = 0.5626639618043182
$ws.Cells.Item(15, 17).Value = 6.360287167658667
$ws.Cells.Item(15, 18).Value = 57.242584508928
$ws.Cells.Item(15, 19).Value = 0.003331172574550431
$ws.Cells.Item(15, 20).Value = 0.00505503953853784

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 0.279608
$ws.Cells.Item(16, 8).Value = 0.838824
$ws.Cells.Item(16, 9).Value = 0.0069281806301764
$ws.Cells.Item(16, 10).Value = 0.008984118197880725
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.281676
$ws.Cells.Item(16, 14).Value = 0.845028
$ws.Cells.Item(16, 15).Value = 0.005953887765346076
$ws.Cells.Item(16, 16).Value = 0.006967417149436334
$ws.Cells.Item(16, 17).Value = 0.078758863008
$ws.Cells.Item(16, 18).Value = 0.708829767072
$ws.Cells.Item(16, 19).Value = [double]"4.124960989011494E-05"
$ws.Cells.Item(16, 20).Value = [double]"6.259609920447722E-05"

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 0.279608
$ws.Cells.Item(17, 8).Value = 0.838824
$ws.Cells.Item(17, 9).Value = 0.0069281806301764
$ws.Cells.Item(17, 10).Value = 0.008984118197880725
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.454891
$ws.Cells.Item(17, 14).Value = 1.364673
$ws.Cells.Item(17, 15).Value = 0.009615196038945605
$ws.Cells.Item(17, 16).Value = 0.01125198935842686
$ws.Cells.Item(17, 17).Value = 0.127191162728
$ws.Cells.Item(17, 18).Value = 1.144720464552
$ws.Cells.Item(17, 19).Value = [double]"6.661581495237178E-05"
$ws.Cells.Item(17, 20).Value = 0.000101089202357403

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 0.279608
$ws.Cells.Item(18, 8).Value = 0.838824
$ws.Cells.Item(18, 9).Value = 0.0069281806301764
$ws.Cells.Item(18, 10).Value = 0.008984118197880725
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 3.179912
$ws.Cells.Item(18, 14).Value = 9.539736
$ws.Cells.Item(18, 15).Value = 0.06721495317910355
$ws.Cells.Item(18, 16).Value = 0.07865694415746599
$ws.Cells.Item(18, 17).Value = 0.889128834496
$ws.Cells.Item(18, 18).Value = 8.002159510463999
$ws.Cells.Item(18, 19).Value = 0.0004656773366736788
$ws.Cells.Item(18, 20).Value = 0.0007066632833947782

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 0.279608
$ws.Cells.Item(19, 8).Value = 0.838824
$ws.Cells.Item(19, 9).Value = 0.0069281806301764
$ws.Cells.Item(19, 10).Value = 0.008984118197880725
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 13).Value = 18.0933415
$ws.Cells.Item(19, 14).Value = 36.186683
$ws.Cells.Item(19, 15).Value = 0.3824455210634858
$ws.Cells.Item(19, 16).Value = 0.2983661082418763
$ws.Cells.Item(19, 17).Value = 5.059043030132001
$ws.Cells.Item(19, 18).Value = 30.354258180792
$ws.Cells.Item(19, 19).Value = 0.002649651651129763
$ws.Cells.Item(19, 20).Value = 0.002680556382686691

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 0.4529403333333333
$ws.Cells.Item(20, 8).Value = 1.358821
$ws.Cells.Item(20, 9).Value = 0.0112230424166177
$ws.Cells.Item(20, 10).Value = 0.01455348019818518
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 13).Value = 2.552614
$ws.Cells.Item(20, 14).Value = 5.105228
$ws.Cells.Item(20, 15).Value = 0.05395552785558979
$ws.Cells.Item(20, 16).Value = 0.04209357928847631
$ws.Cells.Item(20, 17).Value = 1.156181836031333
$ws.Cells.Item(20, 18).Value = 6.937091016188001
$ws.Cells.Item(20, 19).Value = 0.000605545177734282
$ws.Cells.Item(20, 20).Value = 0.0006126080726455778

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 0.4529403333333333
$ws.Cells.Item(21, 8).Value = 1.358821
$ws.Cells.Item(21, 9).Value = 0.0112230424166177
$ws.Cells.Item(21, 10).Value = 0.01455348019818518
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 22.74715733333333
$ws.Cells.Item(21, 14).Value = 68.241472
$ws.Cells.Item(21, 15).Value = 0.4808149140975291
$ws.Cells.Item(21, 16).Value = 0.5626639618043182
$ws.Cells.Item(21, 17).Value = 10.30310502494578
$ws.Cells.Item(21, 18).Value = 92.72794522451201
$ws.Cells.Item(21, 19).Value = 0.005396206175458965
$ws.Cells.Item(21, 20).Value = 0.008188718826351566

$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 0.4529403333333333
$ws.Cells.Item(22, 8).Value = 1.358821
$ws.Cells.Item(22, 9).Value = 0.0112230424166177
$ws.Cells.Item(22, 10).Value = 0.01455348019818518
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 0.281676
$ws.Cells.Item(22, 14).Value = 0.845028
$ws.Cells.Item(22, 15).Value = 0.005953887765346076
$ws.Cells.Item(22, 16).Value = 0.006967417149436334
$ws.Cells.Item(22, 17).Value = 0.127582421332
$ws.Cells.Item(22, 18).Value = 1.148241791988
$ws.Cells.Item(22, 19).Value = [double]"6.682073493426018E-05"
$ws.Cells.Item(22, 20).Value = 0.0001014001675168175

$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 0.4529403333333333
$ws.Cells.Item(23, 8).Value = 1.358821
$ws.Cells.Item(23, 9).Value = 0.0112230424166177
$ws.Cells.Item(23, 10).Value = 0.01455348019818518
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 0.454891
$ws.Cells.Item(23, 14).Value = 1.364673
$ws.Cells.Item(23, 15).Value = 0.009615196038945605
$ws.Cells.Item(23, 16).Value = 0.01125198935842686
$ws.Cells.Item(23, 17).Value = 0.2060384811703333
$ws.Cells.Item(23, 18).Value = 1.854346330533
$ws.Cells.Item(23, 19).Value = 0.000107911752989181
$ws.Cells.Item(23, 20).Value = 0.0001637556043180556

$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 0.4529403333333333
$ws.Cells.Item(24, 8).Value = 1.358821
$ws.Cells.Item(24, 9).Value = 0.0112230424166177
$ws.Cells.Item(24, 10).Value = 0.01455348019818518
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 3.179912
$ws.Cells.Item(24, 14).Value = 9.539736
$ws.Cells.Item(24, 15).Value = 0.06721495317910355
$ws.Cells.Item(24, 16).Value = 0.07865694415746599
$ws.Cells.Item(24, 17).Value = 1.440310401250667
$ws.Cells.Item(24, 18).Value = 12.962793611256
$ws.Cells.Item(24, 19).Value = 0.0007543562705600518
$ws.Cells.Item(24, 20).Value = 0.001144732279245439

$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 0.4529403333333333
$ws.Cells.Item(25, 8).Value = 1.358821
$ws.Cells.Item(25, 9).Value = 0.0112230424166177
$ws.Cells.Item(25, 10).Value = 0.01455348019818518
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 13).Value = 18.0933415
$ws.Cells.Item(25, 14).Value = 36.186683
$ws.Cells.Item(25, 15).Value = 0.3824455210634858
$ws.Cells.Item(25, 16).Value = 0.2983661082418763
$ws.Cells.Item(25, 17).Value = 8.195204130123834
$ws.Cells.Item(25, 18).Value = 49.171224780743
$ws.Cells.Item(25, 19).Value = 0.00429220230494096
$ws.Cells.Item(25, 20).Value = 0.004342265248107722

$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 1.651045666666667
$ws.Cells.Item(26, 8).Value = 4.953137
$ws.Cells.Item(26, 9).Value = 0.04090992606555134
$ws.Cells.Item(26, 10).Value = 0.05304994642296398
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 2.552614
$ws.Cells.Item(26, 14).Value = 5.105228
$ws.Cells.Item(26, 15).Value = 0.05395552785558979
$ws.Cells.Item(26, 16).Value = 0.04209357928847631
$ws.Cells.Item(26, 17).Value = 4.214482283372667
$ws.Cells.Item(26, 18).Value = 25.286893700236
$ws.Cells.Item(26, 19).Value = 0.002207316655399974
$ws.Cells.Item(26, 20).Value = 0.002233062126004454

$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 7).Value = 1.651045666666667
$ws.Cells.Item(27, 8).Value = 4.953137
$ws.Cells.Item(27, 9).Value = 0.04090992606555134
$ws.Cells.Item(27, 10).Value = 0.05304994642296398
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 13).Value = 22.74715733333333
$ws.Cells.Item(27, 14).Value = 68.241472
$ws.Cells.Item(27, 15).Value = 0.4808149140975291
$ws.Cells.Item(27, 16).Value = 0.5626639618043182
$ws.Cells.Item(27, 17).Value = 37.55659554418489
$ws.Cells.Item(27, 18).Value = 338.009359897664
$ws.Cells.Item(27, 19).Value = 0.01967010258694434
$ws.Cells.Item(27, 20).Value = 0.02984929302785173

$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 7).Value = 1.651045666666667
$ws.Cells.Item(28, 8).Value = 4.953137
$ws.Cells.Item(28, 9).Value = 0.04090992606555134
$ws.Cells.Item(28, 10).Value = 0.05304994642296398
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 13).Value = 0.281676
$ws.Cells.Item(28, 14).Value = 0.845028
$ws.Cells.Item(28, 15).Value = 0.005953887765346076
$ws.Cells.Item(28, 16).Value = 0.006967417149436334
$ws.Cells.Item(28, 17).Value = 0.465059939204
$ws.Cells.Item(28, 18).Value = 4.185539452836
$ws.Cells.Item(28, 19).Value = 0.0002435731082828987
$ws.Cells.Item(28, 20).Value = 0.0003696211064840379

$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 7).Value = 1.651045666666667
$ws.Cells.Item(29, 8).Value = 4.953137
$ws.Cells.Item(29, 9).Value = 0.04090992606555134
$ws.Cells.Item(29, 10).Value = 0.05304994642296398
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 13).Value = 0.454891
$ws.Cells.Item(29, 14).Value = 1.364673
$ws.Cells.Item(29, 15).Value = 0.009615196038945605
$ws.Cells.Item(29, 16).Value = 0.01125198935842686
$ws.Cells.Item(29, 17).Value = 0.7510458143556666
$ws.Cells.Item(29, 18).Value = 6.759412329201
$ws.Cells.Item(29, 19).Value = 0.0003933569590590468
$ws.Cells.Item(29, 20).Value = 0.0005969174326163056

$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 7).Value = 1.651045666666667
$ws.Cells.Item(30, 8).Value = 4.953137
$ws.Cells.Item(30, 9).Value = 0.04090992606555134
$ws.Cells.Item(30, 10).Value = 0.05304994642296398
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 13).Value = 3.179912
$ws.Cells.Item(30, 14).Value = 9.539736
$ws.Cells.Item(30, 15).Value = 0.06721495317910355
$ws.Cells.Item(30, 16).Value = 0.07865694415746599
$ws.Cells.Item(30, 17).Value = 5.250179927981333
$ws.Cells.Item(30, 18).Value = 47.251619351832
$ws.Cells.Item(30, 19).Value = 0.002749758765056621
$ws.Cells.Item(30, 20).Value = 0.00417274667334764

$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 7).Value = 1.651045666666667
$ws.Cells.Item(31, 8).Value = 4.953137
$ws.Cells.Item(31, 9).Value = 0.04090992606555134
$ws.Cells.Item(31, 10).Value = 0.05304994642296398
$ws.Cells.Item(31, 11).Value = 2
$ws.Cells.Item(31, 13).Value = 18.0933415
$ws.Cells.Item(31, 14).Value = 36.186683
$ws.Cells.Item(31, 15).Value = 0.3824455210634858
$ws.Cells.Item(31, 16).Value = 0.2983661082418763
$ws.Cells.Item(31, 17).Value = 29.87293307909517
$ws.Cells.Item(31, 18).Value = 179.237598474571
$ws.Cells.Item(31, 19).Value = 0.01564581799080846
$ws.Cells.Item(31, 20).Value = 0.01582830605665981

$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 7).Value = 25.8827875
$ws.Cells.Item(32, 8).Value = 51.765575
$ws.Cells.Item(32, 9).Value = 0.6413286708980853
$ws.Cells.Item(32, 10).Value = 0.5544286338746381
$ws.Cells.Item(32, 11).Value = 2
$ws.Cells.Item(32, 13).Value = 2.552614
$ws.Cells.Item(32, 14).Value = 5.105228
$ws.Cells.Item(32, 15).Value = 0.05395552785558979
$ws.Cells.Item(32, 16).Value = 0.04209357928847631
$ws.Cells.Item(32, 17).Value = 66.068765731525
$ws.Cells.Item(32, 18).Value = 264.2750629261
$ws.Cells.Item(32, 19).Value = 0.03460322696723002
$ws.Cells.Item(32, 20).Value = 0.02333788565980368

$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 7).Value = 25.8827875
$ws.Cells.Item(33, 8).Value = 51.765575
$ws.Cells.Item(33, 9).Value = 0.6413286708980853
$ws.Cells.Item(33, 10).Value = 0.5544286338746381
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 13).Value = 22.74715733333333
$ws.Cells.Item(33, 14).Value = 68.241472
$ws.Cells.Item(33, 15).Value = 0.4808149140975291
$ws.Cells.Item(33, 16).Value = 0.5626639618043182
$ws.Cells.Item(33, 17).Value = 588.7598394877333
$ws.Cells.Item(33, 18).Value = 3532.5590369264
$ws.Cells.Item(33, 19).Value = 0.3083603898061454
$ws.Cells.Item(33, 20).Value = 0.3119570116736597

$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(34, 7).Value = 25.8827875
$ws.Cells.Item(34, 8).Value = 51.765575
$ws.Cells.Item(34, 9).Value = 0.6413286708980853
$ws.Cells.Item(34, 10).Value = 0.5544286338746381
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 13).Value = 0.281676
$ws.Cells.Item(34, 14).Value = 0.845028
$ws.Cells.Item(34, 15).Value = 0.005953887765346076
$ws.Cells.Item(34, 16).Value = 0.006967417149436334
$ws.Cells.Item(34, 17).Value = 7.290560051849999
$ws.Cells.Item(34, 18).Value = 43.7433603111
$ws.Cells.Item(34, 19).Value = 0.00381839892722577
$ws.Cells.Item(34, 20).Value = 0.003862935571796712

$ws.Cells.Item(35, 5).Value = 2
$ws.Cells.Item(35, 7).Value = 25.8827875
$ws.Cells.Item(35, 8).Value = 51.765575
$ws.Cells.Item(35, 9).Value = 0.6413286708980853
$ws.Cells.Item(35, 10).Value = 0.5544286338746381
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 13).Value = 0.454891
$ws.Cells.Item(35, 14).Value = 1.364673
$ws.Cells.Item(35, 15).Value = 0.009615196038945605
$ws.Cells.Item(35, 16).Value = 0.01125198935842686
$ws.Cells.Item(35, 17).Value = 11.7738470886625
$ws.Cells.Item(35, 18).Value = 70.643082531975
$ws.Cells.Item(35, 19).Value = 0.006166500896081519
$ws.Cells.Item(35, 20).Value = 0.006238425088364569

$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 7).Value = 25.8827875
$ws.Cells.Item(36, 8).Value = 51.765575
$ws.Cells.Item(36, 9).Value = 0.6413286708980853
$ws.Cells.Item(36, 10).Value = 0.5544286338746381
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 13).Value = 3.179912
$ws.Cells.Item(36, 14).Value = 9.539736
$ws.Cells.Item(36, 15).Value = 0.06721495317910355
$ws.Cells.Item(36, 16).Value = 0.07865694415746599
$ws.Cells.Item(36, 17).Value = 82.30498656469999
$ws.Cells.Item(36, 18).Value = 493.8299193882
$ws.Cells.Item(36, 19).Value = 0.04310687658683151
$ws.Cells.Item(36, 20).Value = 0.04360966209397757

$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 7).Value = 25.8827875
$ws.Cells.Item(37, 8).Value = 51.765575
$ws.Cells.Item(37, 9).Value = 0.6413286708980853
$ws.Cells.Item(37, 10).Value = 0.5544286338746381
$ws.Cells.Item(37, 11).Value = 2
$ws.Cells.Item(37, 13).Value = 18.0933415
$ws.Cells.Item(37, 14).Value = 36.186683
$ws.Cells.Item(37, 15).Value = 0.3824455210634858
$ws.Cells.Item(37, 16).Value = 0.2983661082418763
$ws.Cells.Item(37, 17).Value = 468.3061132094313
$ws.Cells.Item(37, 18).Value = 1873.224452837725
$ws.Cells.Item(37, 19).Value = 0.2452732777145711
$ws.Cells.Item(37, 20).Value = 0.1654227137870359
